$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting existing data (D:K) right to (E:L)
$ws.Columns("D").Insert()

# Copy cell formatting (number format/font/style) from the old D column (now E) into the new D column
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new column D with the latest reporting period values
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 798100
$ws.Range("D9").Value = 489500
$ws.Range("D10").Value = 308700
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 157500
$ws.Range("D15").Value = 58700
$ws.Range("D17").Value = 878500
$ws.Range("D18").Value = -80400
$ws.Range("D20").Value = 62900
$ws.Range("D21").Value = 41200
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = -17500
$ws.Range("D24").Value = 800
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = -18300
$ws.Range("D27").Value = -64800
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -62900
$ws.Range("D33").Value = -64800
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -64800
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 931800
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = "NA"
$ws.Range("D44").Value = 22600
$ws.Range("D45").Value = "NA"
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 1339400
$ws.Range("D48").Value = 2373100
$ws.Range("D49").Value = 162500
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 43700
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 5014300
$ws.Range("D57").Value = 0
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 162200
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 3609100
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 4151300
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -2472100
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 863000
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = -64800
$ws.Range("D83").Value = 58700
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = -24100
$ws.Range("D91").Value = -208500
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = 778900
$ws.Range("D96").Value = -44700
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -457900
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 296800

# A handful of "Total"/"Change" rows were recalculated for the 3 most-recent periods
# (their D/E/F values are not simple carries of the prior period data)
$ws.Range("E89").Value = 101500
$ws.Range("F89").Value = 29500
$ws.Range("E94").Value = 263100
$ws.Range("F94").Value = 465000
$ws.Range("E100").Value = -41500
$ws.Range("F100").Value = -877700
$ws.Range("E102").Value = 323100
$ws.Range("F102").Value = -383100
